$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.775841999999999
$ws.Range("H2").Value = 8.327525999999999
$ws.Range("I2").Value = 0.0624750527258915
$ws.Range("J2").Value = 0.0624750527258915
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 169.4418546075327
$ws.Range("R2").Value = 1524.976691467794
$ws.Range("S2").Value = 0.01276748587007356
$ws.Range("T2").Value = 0.01276748587007356

$ws.Range("G3").Value = 2.775841999999999
$ws.Range("H3").Value = 8.327525999999999
$ws.Range("I3").Value = 0.0624750527258915
$ws.Range("J3").Value = 0.0624750527258915
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 295.1121599303719
$ws.Range("R3").Value = 2656.009439373348
$ws.Range("S3").Value = 0.02223677461938268
$ws.Range("T3").Value = 0.02223677461938268

$ws.Range("G4").Value = 2.775841999999999
$ws.Range("H4").Value = 8.327525999999999
$ws.Range("I4").Value = 0.0624750527258915
$ws.Range("J4").Value = 0.0624750527258915
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 364.5746728406585
$ws.Range("R4").Value = 3281.172055565927
$ws.Range("S4").Value = 0.02747079223643525
$ws.Range("T4").Value = 0.02747079223643526

$ws.Range("I5").Value = 0.2652892219050753
$ws.Range("J5").Value = 0.2652892219050753
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 719.5047591909622
$ws.Range("R5").Value = 6475.54283271866
$ws.Range("S5").Value = 0.05421486248305564
$ws.Range("T5").Value = 0.05421486248305564

$ws.Range("I6").Value = 0.2652892219050753
$ws.Range("J6").Value = 0.2652892219050753
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.09442451633193684
$ws.Range("T6").Value = 0.09442451633193685

$ws.Range("I7").Value = 0.2652892219050753
$ws.Range("J7").Value = 0.2652892219050753
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 1548.10163520059
$ws.Range("R7").Value = 13932.91471680531
$ws.Range("S7").Value = 0.1166498430900828
$ws.Range("T7").Value = 0.1166498430900828

$ws.Range("G8").Value = 29.86824466666667
$ws.Range("H8").Value = 89.60473400000001
$ws.Range("I8").Value = 0.6722357253690333
$ws.Range("J8").Value = 0.6722357253690333
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 1823.20563281035
$ws.Range("R8").Value = 16408.85069529315
$ws.Range("S8").Value = 0.1373789977043242
$ws.Range("T8").Value = 0.1373789977043242

$ws.Range("G9").Value = 29.86824466666667
$ws.Range("H9").Value = 89.60473400000001
$ws.Range("I9").Value = 0.6722357253690333
$ws.Range("J9").Value = 0.6722357253690333
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 3175.426482094015
$ws.Range("R9").Value = 28578.83833884613
$ws.Range("S9").Value = 0.2392691748771168
$ws.Range("T9").Value = 0.2392691748771168

$ws.Range("G10").Value = 29.86824466666667
$ws.Range("H10").Value = 89.60473400000001
$ws.Range("I10").Value = 0.6722357253690333
$ws.Range("J10").Value = 0.6722357253690333
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 3922.847744098817
$ws.Range("R10").Value = 35305.62969688935
$ws.Range("S10").Value = 0.2955875527875922
$ws.Range("T10").Value = 0.2955875527875923
